# Update "想去人数" (want-to-go count) values in column F
# on the "展览" (Exhibition) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7
$ws1.Range("F7").Value = 7661
$ws1.Range("F9").Value = 203
$ws1.Range("F10").Value = 1079
$ws1.Range("F11").Value = 654
$ws1.Range("F12").Value = 12
$ws1.Range("F14").Value = 173
$ws1.Range("F15").Value = 5
$ws1.Range("F16").Value = 202
$ws1.Range("F17").Value = 750

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7
$ws4.Range("F8").Value = 7661
$ws4.Range("F10").Value = 203
$ws4.Range("F11").Value = 1079
$ws4.Range("F12").Value = 654
$ws4.Range("F13").Value = 12
$ws4.Range("F15").Value = 173
$ws4.Range("F16").Value = 5
$ws4.Range("F17").Value = 202
$ws4.Range("F18").Value = 750
